$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff swaps the data between row 2 and row 3 for the date column (D)
# and the price-related columns (M, N, O, P, S), leaving everything else
# (A, B, C, E, F, G, H, I, J, K, L, Q, R, T) unchanged.

$ws.Range("D2").Value = 44322
$ws.Range("D3").Value = 44365

$ws.Range("M2").Value = 600
$ws.Range("N2").Value = 1500
$ws.Range("O2").Value = 1600
$ws.Range("P2").Value = 1550
$ws.Range("S2").Value = 1550

$ws.Range("M3").Value = 900
$ws.Range("N3").Value = 1200
$ws.Range("O3").Value = 1400
$ws.Range("P3").Value = 1300
$ws.Range("S3").Value = 1300
